{"js": "// Apply the \"Added many more features\" edit: refresh the title, the\n// three \"What we like\" / \"What we don't like\" bullet lists, the bolded\n// title repeated near the end, and the closing italic summary blurb.\n\nconst replacements = [\n  {\n    find: \"Play Bat Stax for Free - A Halloween-Themed Slot Game\",\n    replace: \"Play Bat Stax Free - Exciting Halloween-themed Slot Game\",\n  },\n  {\n    find: \"Excellent graphics creating a chilling atmosphere\",\n    replace: \"Great graphics and chilling atmosphere\",\n  },\n  {\n    find: \"Medium volatility level ensuring good winnings\",\n    replace: \"Medium volatility ensures frequent winnings\",\n  },\n  {\n    find: \"Free Spins feature with potential for grid to be filled with symbols\",\n    replace: \"Exciting bonus features, including free spins\",\n  },\n  {\n    find: \"Accessible betting range from 1 token per line up to 200 tokens\",\n    replace: \"Accessible betting range for all players\",\n  },\n  {\n    find: \"Autoplay feature with no predefined number of spins\",\n    replace: \"Autoplay feature lacks a predefined number of spins\",\n  },\n  {\n    find: \"Not the highest RTP level at 95.2%\",\n    replace: \"Limited paylines for more advanced players\",\n  },\n  {\n    find: \"Read our review of Bat Stax, a Halloween-themed slot game with free spin feature. Play for free and enjoy the eerie music and scary animations.\",\n    replace: \"Play Bat Stax for free and experience the thrilling Halloween atmosphere with exciting bonus features.\",\n  },\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Added many more features\" edit: refresh the title, the\n# \"What we like\" / \"What we don't like\" bullet lists, the bolded title\n# repeated near the end, and the closing italic summary blurb.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"Play Bat Stax for Free - A Halloween-Themed Slot Game\"; Replace = \"Play Bat Stax Free - Exciting Halloween-themed Slot Game\" },\n    @{ Find = \"Excellent graphics creating a chilling atmosphere\"; Replace = \"Great graphics and chilling atmosphere\" },\n    @{ Find = \"Medium volatility level ensuring good winnings\"; Replace = \"Medium volatility ensures frequent winnings\" },\n    @{ Find = \"Free Spins feature with potential for grid to be filled with symbols\"; Replace = \"Exciting bonus features, including free spins\" },\n    @{ Find = \"Accessible betting range from 1 token per line up to 200 tokens\"; Replace = \"Accessible betting range for all players\" },\n    @{ Find = \"Autoplay feature with no predefined number of spins\"; Replace = \"Autoplay feature lacks a predefined number of spins\" },\n    @{ Find = \"Not the highest RTP level at 95.2%\"; Replace = \"Limited paylines for more advanced players\" },\n    @{ Find = \"Read our review of Bat Stax, a Halloween-themed slot game with free spin feature. Play for free and enjoy the eerie music and scary animations.\"; Replace = \"Play Bat Stax for free and experience the thrilling Halloween atmosphere with exciting bonus features.\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.Text = $r.Find\n    $find.Replacement.Text = $r.Replace\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
